$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally listed years 1999-2016 (rows 2-19). The data was
# "rows transposed to columns" — i.e. re-based so the table now starts at
# 2008 (dropping 1999-2007), shifting the 2008-2016 rows up into rows 2-10.
# Deleting rows 2-10 (1999-2007) shifts rows 11-19 (2008-2016) up to 2-10.
$ws.Range("A2:B10").EntireRow.Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Match the new selection recorded in the sheet view.
$ws.Range("D8").Select()
